$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its original text formatting so that
# numeric-looking strings (with thousand separators as dots, fixed decimal
# places, etc.) are not reinterpreted/rounded as floating point numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.748.11'
$ws.Range("E2").Value = '  +1.06%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.246.52'
$ws.Range("E3").Value = '  +0.34%  '

$ws.Range("E4").Value = '  +0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '320.79'
$ws.Range("E5").Value = '  +0.17%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '101.94'
$ws.Range("E6").Value = '  +0.41%  '

$ws.Range("E7").Value = '  -1.13%  '

$ws.Range("E8").Value = '  +0.05%  '

$ws.Range("E9").Value = '  -1.34%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.01'
$ws.Range("E10").Value = '  -0.25%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0833'
$ws.Range("E11").Value = '  +0.52%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.73'
$ws.Range("E12").Value = '  +0.02%  '

$ws.Range("E13").Value = '  -2.32%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.590.64'
$ws.Range("E14").Value = '  +0.29%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.855'
$ws.Range("E15").Value = '  -1.40%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.21'
$ws.Range("E16").Value = '  -1.48%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.246.16'
$ws.Range("E17").Value = '  +0.82%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.705.06'
$ws.Range("E18").Value = '  +1.09%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.56'
$ws.Range("E19").Value = '  -5.81%  '

$ws.Range("E20").Value = '  +2.12%  '

$ws.Range("E21").Value = '  +0.17%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.48'
$ws.Range("E22").Value = '  -0.12%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.16'
$ws.Range("E23").Value = '  -1.11%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '235.72'
$ws.Range("E24").Value = '  -0.85%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.13'
$ws.Range("E25").Value = '  -1.25%  '

$ws.Range("E26").Value = '  +0.14%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.16'
$ws.Range("E27").Value = '  +0.94%  '

$ws.Range("E28").Value = '  -2.85%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '37.12'
$ws.Range("E29").Value = '  +3.74%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.29'
$ws.Range("E30").Value = '  -1.81%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '160.06'
$ws.Range("E31").Value = '  +4.21%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.16'
$ws.Range("E32").Value = '  -1.44%  '

$ws.Range("E33").Value = '  -3.10%  '

$ws.Range("E34").Value = '  -1.49%  '

$ws.Range("E35").Value = '  +9.71%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.07'
$ws.Range("E36").Value = '  -3.37%  '

$ws.Range("E37").Value = '  -1.04%  '

$ws.Range("E38").Value = '  -3.01%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.76'
$ws.Range("E39").Value = '  +0.82%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.24'
$ws.Range("E40").Value = '  -4.88%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '15.81'
$ws.Range("E41").Value = '  +19.80%  '

$ws.Range("E42").Value = '  -2.21%  '

$ws.Range("E43").Value = '  +0.12%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.806.40'
$ws.Range("E44").Value = '  +1.47%  '

$ws.Range("B45").Value = 'Algorand'
$ws.Range("C45").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.200'
$ws.Range("E45").Value = '  -2.94%  '

$ws.Range("B46").Value = 'ordi'
$ws.Range("C46").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '75.91'
$ws.Range("E46").Value = '  -0.43%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '82.25'
$ws.Range("E47").Value = '  -5.00%  '

$ws.Range("E48").Value = '  -2.42%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '58.59'
$ws.Range("E49").Value = '  -1.09%  '

$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '103.60'
$ws.Range("E50").Value = '  -0.12%  '

$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.67'
$ws.Range("E51").Value = '  +4.91%  '

